$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ROW3
$ws.Range("G3").Value = 1.8
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 2.5
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 5
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 7.5
$ws.Range("Z3").Value = 13
$ws.Range("AA3").Value = 17
$ws.Range("AC3").Value = 7.5
$ws.Range("AD3").Value = 6.5
$ws.Range("AE3").Value = 19
$ws.Range("AG3").Value = 10
$ws.Range("AM3").Value = 1250
$ws.Range("AN3").Value = 3.6
$ws.Range("AO3").Value = 10
$ws.Range("AQ3").Value = 34
$ws.Range("AT3").Value = 2.5
$ws.Range("AX3").Value = 26
$ws.Range("AZ3").Value = 101
# ROW8
$ws.Range("G8").Value = 3.65
$ws.Range("I8").Value = 1.9
$ws.Range("J8").Value = 4.05
$ws.Range("K8").Value = 2.12
$ws.Range("L8").Value = 2.47
$ws.Range("O8").Value = 1.24
$ws.Range("P8").Value = 3.3
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 2.02
$ws.Range("W8").Value = 12
$ws.Range("X8").Value = 21
$ws.Range("Y8").Value = 12.5
$ws.Range("Z8").Value = 55
$ws.Range("AA8").Value = 32
$ws.Range("AB8").Value = 35
$ws.Range("AC8").Value = 11.25
$ws.Range("AE8").Value = 13
$ws.Range("AH8").Value = 9.75
$ws.Range("AJ8").Value = 17
$ws.Range("AK8").Value = 14.5
$ws.Range("AN8").Value = 5.5
$ws.Range("AO8").Value = 20
$ws.Range("AP8").Value = 25
$ws.Range("AQ8").Value = 100
$ws.Range("AR8").Value = 120
$ws.Range("AS8").Value = 300
$ws.Range("AW8").Value = 3.8
$ws.Range("AX8").Value = 9.5
$ws.Range("AY8").Value = 17.5
$ws.Range("AZ8").Value = 35
$ws.Range("BA8").Value = 60
# ROW11
$ws.Range("G11").Value = 1.42
$ws.Range("H11").Value = 4.2
$ws.Range("I11").Value = 6.25
$ws.Range("J11").Value = 2
$ws.Range("K11").Value = 2.25
$ws.Range("L11").Value = 7
$ws.Range("U11").Value = 2.2
$ws.Range("V11").Value = 1.62
$ws.Range("Y11").Value = 8.5
$ws.Range("Z11").Value = 9.5
$ws.Range("AD11").Value = 8.5
$ws.Range("AU11").Value = 9.5
$ws.Range("AW11").Value = 8
# ROW23
$ws.Range("O23").Value = 1.14
$ws.Range("P23").Value = 5.5
$ws.Range("BD23").Value = 151
# ROW25
$ws.Range("I25").Value = 2.77
$ws.Range("T25").Value = 3.25
$ws.Range("X25").Value = 13
$ws.Range("AT25").Value = 3.25
# ROW26
$ws.Range("G26").Value = 2.5
$ws.Range("H26").Value = 2.7
$ws.Range("L26").Value = 3.85
$ws.Range("S26").Value = 1.53
$ws.Range("T26").Value = 2.2
$ws.Range("AC26").Value = 6.1
$ws.Range("AH26").Value = 15
$ws.Range("AT26").Value = 2.18
$ws.Range("AU26").Value = 7.2
$ws.Range("AW26").Value = 4.85
$ws.Range("AX26").Value = 19
$ws.Range("AY26").Value = 29
$ws.Range("BA26").Value = 175
# ROW27
$ws.Range("H27").Value = 2.7
$ws.Range("I27").Value = 3.55
$ws.Range("J27").Value = 2.95
$ws.Range("L27").Value = 4.1
$ws.Range("S27").Value = 1.52
$ws.Range("T27").Value = 2.22
$ws.Range("U27").Value = 1.91
$ws.Range("V27").Value = 1.7
$ws.Range("W27").Value = 6
$ws.Range("Y27").Value = 9.25
$ws.Range("Z27").Value = 24
$ws.Range("AD27").Value = 5.4
$ws.Range("AE27").Value = 15.5
$ws.Range("AG27").Value = 8.25
$ws.Range("AH27").Value = 18
$ws.Range("AI27").Value = 12
$ws.Range("AJ27").Value = 55
$ws.Range("AP27").Value = 22
$ws.Range("AR27").Value = 100
$ws.Range("AS27").Value = 350
$ws.Range("AT27").Value = 2.2
$ws.Range("AU27").Value = 7
$ws.Range("AV27").Value = 70
$ws.Range("AW27").Value = 5.2
$ws.Range("AX27").Value = 21
$ws.Range("AY27").Value = 28
$ws.Range("BA27").Value = 150
